$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("DR. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $rng.Start
Write-Host "start=$start text=[$($rng.Text)]"

# Slices (character offsets, determined via probing):
#  br+D    : start-1 .. start+1
#  r       : start+1 .. start+2
#  ". "    : start+2 .. start+4
#  <mname> : start+4 .. start+11
#  " "     : start+11 .. start+12
#  <mfname>: start+12 .. start+20
#  " "     : start+20 .. start+21
#  <msname>: start+21 .. start+29
$slices = @(
  @($start - 1, $start + 1),
  @($start + 1, $start + 2),
  @($start + 2, $start + 4),
  @($start + 4, $start + 11),
  @($start + 11, $start + 12),
  @($start + 12, $start + 20),
  @($start + 20, $start + 21),
  @($start + 21, $start + 29)
)

foreach ($s in $slices) {
  $r = $d.Range($s[0], $s[1])
  Write-Host "slice [$($s[0]),$($s[1])) = [$($r.Text)]"
  $r.Font.Bold = 1
}

$r2b = $d.Range($start+1, $start+2)
$res = $r2b.Find.Execute("R", $true, $false, $false, $false, $false, $true, 0, $false, "r", 1)
Write-Host "replace result=$res"

foreach ($s in $slices) {
  $r = $d.Range($s[0], $s[1])
  $r.Font.Bold = 0
}

Write-Host "Final text: [$($d.Range($start-1, $start+29).Text)]"
